# This cell was misclassified as a "Context" term but is actually "Neither".
# Update the raw data sheet: one unit moves from Context (col B) to Neither (col C)
# for the first contributor row (row 2). Downstream formulas on the stats sheet
# and the chart will recalculate automatically.

$wb = $excel.ActiveWorkbook

$rawSheet = $wb.Worksheets.Item("raw")
$statsSheet = $wb.Worksheets.Item("ugrad-009-01-stats-20")

$rawSheet.Range("B2").Value = 19
$rawSheet.Range("C2").Value = 41

# Force recalculation so dependent formulas / chart caches refresh.
$excel.CalculateFullRebuild()

# Update selections to reflect where the edit was made, and switch the
# active sheet to "raw" (matching the workbook's activeTab/tabSelected change).
$statsSheet.Range("B2:D2").Select()
$rawSheet.Activate()
$rawSheet.Range("A2:C2").Select()
